# "Delete Meta & NoAttri from templates"
#
# The "Label" sheet stores, per Kategorie (column group), a header in row 3
# plus the allowed Merkmal values beneath it. Two whole categories - "Meta"
# (column T) and "NoAttribute" (column U) - are being retired: their header
# cells are cleared and all their data rows are cleared too. Since they are
# also entries in the alphabetically sorted "Attribute" list in column B
# (B1:B16), those two rows (B10:B11, "Meta"/"NoAttribute") are removed and
# the rest of that list shifts up, shrinking the list (and the matching
# "Attribute" named range) from B1:B16 down to B1:B14.
#
# The two corresponding defined names ("Meta", "NoAttribute") are removed
# entirely, and the "Attribute" named range is shrunk to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Label")

# 1. Clear out the "Meta" (T) and "NoAttribute" (U) columns' header + data.
#    Row 3 holds the header text (style s=5 stays on the now-empty cell);
#    rows 4-9 hold the data under "Meta", row 4 alone holds data under
#    "NoAttribute" - ClearContents leaves the cell itself (and any style)
#    in place but empties the value, matching the target exactly for T3/U3
#    while T4:T9/U4 (default style) disappear outright once truly empty.
$ws.Range("T3:U9").ClearContents()

# 2. Remove the "Meta" / "NoAttribute" rows from the Attribute list and
#    shift the remaining entries (Nutzung_Widmung ... Vorbauten) up.
#    (A plain ranged Delete/xlShiftUp here would drag the unrelated H/I/P
#    columns of those rows along with it, so the shift is done by hand,
#    column B only, then the vacated B15:B16 tail cells are cleared away.)
for ($i = 12; $i -le 16; $i++) {
    $ws.Range("B" + ($i - 2)).Value2 = $ws.Range("B$i").Value2
}
$ws.Range("B15:B16").Clear() | Out-Null

# 3. Defined names: drop "Meta" and "NoAttribute", and shrink "Attribute".
$wb.Names.Item("Meta").Delete()
$wb.Names.Item("NoAttribute").Delete()
$wb.Names.Item("Attribute").RefersTo = "=Label!`$B`$1:`$B`$14"

# 4. Cosmetic row-height follow-up that Excel itself performs once the
#    "Meta"/"NoAttribute" columns stop holding the tallest content in row 3
#    and the Attribute list no longer reaches rows 15/16 (now blank again).
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Rows.Item(14).RowHeight = 13.8
$ws.Rows.Item(15).RowHeight = 13.8
$ws.Rows.Item(16).RowHeight = 13.8

# 5. Restore the Label sheet's view to a neutral top-left / selection state
#    (without stealing the active-tab away from Dataset, which stays the
#    workbook's selected sheet both before and after this edit).
$ws.Range("B10").Select() | Out-Null
$wb.Worksheets.Item("Dataset").Activate() | Out-Null
